$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6681570004208766
$ws.Range("C2").Value = 0.1604053655967164
$ws.Range("D2").Value = 0.04564022844408555
$ws.Range("E2").Value = 0.1136851970799135
$ws.Range("F2").Value = 0.9761275492257369
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("K2").Value = 0.3453905516442148
$ws.Range("L2").Value = 0.1949250355620507
$ws.Range("N2").Value = 1.777081668538464
$ws.Range("O2").Value = 3.483441505578071
$ws.Range("B3").Value = 0.6256801764266697
$ws.Range("C3").Value = 0.1602524579789559
$ws.Range("D3").Value = 0.0437770337975536
$ws.Range("E3").Value = 0.1131101741652856
$ws.Range("F3").Value = 0.9753539534409441
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("K3").Value = 0.3087244553701112
$ws.Range("L3").Value = 0.1878130467408994
$ws.Range("N3").Value = 1.795260466065978
$ws.Range("O3").Value = 3.494724656376036
$ws.Range("B4").Value = 0.5998426459453867
$ws.Range("C4").Value = 0.1601689875636119
$ws.Range("D4").Value = 0.04261964192335199
$ws.Range("E4").Value = 0.1128139854476196
$ws.Range("F4").Value = 0.9754090563247644
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("K4").Value = 0.2862491978764012
$ws.Range("L4").Value = 0.1835499510698924
$ws.Range("N4").Value = 1.806992143303564
$ws.Range("O4").Value = 3.503407944328472
$ws.Range("B5").Value = 0.5893754893461676
$ws.Range("C5").Value = 0.1601376117360331
$ws.Range("D5").Value = 0.04214465020790215
$ws.Range("E5").Value = 0.1127076087529169
$ws.Range("F5").Value = 0.9755648813239333
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("K5").Value = 0.2771003453930376
$ws.Range("L5").Value = 0.1818388534496052
$ws.Range("N5").Value = 1.811916181421156
$ws.Range("O5").Value = 3.507388020435528
$ws.Range("B6").Value = 0.5876411801836241
$ws.Range("C6").Value = 0.1601325617832323
$ws.Range("D6").Value = 0.04206557674898193
$ws.Range("E6").Value = 0.1126908106308555
$ws.Range("F6").Value = 0.97559881335075
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("K6").Value = 0.2755818031805006
$ws.Range("L6").Value = 0.1815563084686715
$ws.Range("N6").Value = 1.812742468205569
$ws.Range("O6").Value = 3.508075583253046
$ws.Range("B7").Value = 0.5997012309481136
$ws.Range("C7").Value = 0.16016855370653
$ws.Range("D7").Value = 0.04261324953695578
$ws.Range("E7").Value = 0.1128124927954275
$ws.Range("F7").Value = 0.975410617727519
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("K7").Value = 0.2861257721453114
$ws.Range("L7").Value = 0.1835267686367104
$ws.Range("N7").Value = 1.807057970453677
$ws.Range("O7").Value = 3.503459832935476
$ws.Range("B8").Value = 0.6534608408236409
$ws.Range("C8").Value = 0.1603504931666215
$ws.Range("D8").Value = 0.04500058810113927
$ws.Range("E8").Value = 0.1134751359153938
$ws.Range("F8").Value = 0.975750833764117
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("K8").Value = 0.3327404938336542
$ws.Range("L8").Value = 0.1924513417285709
$ws.Range("N8").Value = 1.783231451866615
$ws.Range("O8").Value = 3.48696769487924
$ws.Range("B9").Value = 0.7607922397304776
$ws.Range("C9").Value = 0.160789084672011
$ws.Range("D9").Value = 0.0495752740928026
$ws.Range("E9").Value = 0.1152251885736959
$ws.Range("F9").Value = 0.980622176780507
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("K9").Value = 0.4244360343206779
$ws.Range("L9").Value = 0.2107730314314864
$ws.Range("N9").Value = 1.741028368442902
$ws.Range("O9").Value = 3.468550780601646
$ws.Range("B10").Value = 0.8407906717435765
$ws.Range("C10").Value = 0.161160135818065
$ws.Range("D10").Value = 0.05287049434146951
$ws.Range("E10").Value = 0.1167850329405518
$ws.Range("F10").Value = 0.9867636300129163
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("K10").Value = 0.4919627630805508
$ws.Range("L10").Value = 0.2247330074572744
$ws.Range("N10").Value = 1.712775750114286
$ws.Range("O10").Value = 3.4635071857806
$ws.Range("B11").Value = 0.8774279078908478
$ws.Range("C11").Value = 0.1613393189856325
$ws.Range("D11").Value = 0.0543551657374266
$ws.Range("E11").Value = 0.1175540468047345
$ws.Range("F11").Value = 0.9901139800157068
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("K11").Value = 0.522713887515124
$ws.Range("L11").Value = 0.2311919513841616
$ws.Range("N11").Value = 1.700520497415758
$ws.Range("O11").Value = 3.463055735569839
$ws.Range("B12").Value = 0.8913362473992379
$ws.Range("C12").Value = 0.1614086459383088
$ws.Range("D12").Value = 0.05491529280821084
$ws.Range("E12").Value = 0.1178537829382407
$ws.Range("F12").Value = 0.9914626671603912
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("K12").Value = 0.5343628709393045
$ws.Range("L12").Value = 0.2336533325760399
$ws.Range("N12").Value = 1.695965610553354
$ws.Range("O12").Value = 3.463149727143303
$ws.Range("B13").Value = 0.8883393092922915
$ws.Range("C13").Value = 0.1613936498900728
$ws.Range("D13").Value = 0.05479475256212396
$ws.Range("E13").Value = 0.1177888504513263
$ws.Range("F13").Value = 0.9911686468554421
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("K13").Value = 0.5318538755440727
$ws.Range("L13").Value = 0.2331225411532358
$ws.Range("N13").Value = 1.69694276497825
$ws.Range("O13").Value = 3.463117701821034
$ws.Range("B14").Value = 0.8785714656409596
$ws.Range("C14").Value = 0.1613449931260647
$ws.Range("D14").Value = 0.05440128966650093
$ws.Range("E14").Value = 0.1175785354872723
$ws.Range("F14").Value = 0.9902233346552407
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("K14").Value = 0.5236721752428934
$ws.Range("L14").Value = 0.2313941401271364
$ws.Range("N14").Value = 1.700144041971906
$ws.Range("O14").Value = 3.46305815875678
$ws.Range("B15").Value = 0.8725928663614866
$ws.Range("C15").Value = 0.1613153808426802
$ws.Range("D15").Value = 0.05416000995261072
$ws.Range("E15").Value = 0.1174508214333905
$ws.Range("F15").Value = 0.9896547181361655
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("K15").Value = 0.5186611788586504
$ws.Range("L15").Value = 0.230337463276868
$ws.Range("N15").Value = 1.70211610666756
$ws.Range("O15").Value = 3.463056188856058
$ws.Range("B16").Value = 0.8384012237601155
$ws.Range("C16").Value = 0.1611486331463325
$ws.Range("D16").Value = 0.05277317698605799
$ws.Range("E16").Value = 0.1167359703864292
$ws.Range("F16").Value = 0.9865558719802578
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("K16").Value = 0.4899537212524478
$ws.Range("L16").Value = 0.2243130763582855
$ws.Range("N16").Value = 1.713588678272762
$ws.Range("O16").Value = 3.463573767836522
$ws.Range("B17").Value = 0.8174881541654031
$ws.Range("C17").Value = 0.1610489857897406
$ws.Range("D17").Value = 0.05191871009630944
$ws.Range("E17").Value = 0.1163126410417661
$ws.Range("F17").Value = 0.9847973448937779
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("K17").Value = 0.4723506711355867
$ws.Range("L17").Value = 0.2206450378356379
$ws.Range("N17").Value = 1.720779680501834
$ws.Range("O17").Value = 3.464363254555025
$ws.Range("B18").Value = 0.8054826634902099
$ws.Range("C18").Value = 0.1609926504160484
$ws.Range("D18").Value = 0.05142589492628957
$ws.Range("E18").Value = 0.1160747484478186
$ws.Range("F18").Value = 0.9838382769862122
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("K18").Value = 0.4622289883929511
$ws.Range("L18").Value = 0.21854549464706
$ws.Range("N18").Value = 1.724971943806419
$ws.Range("O18").Value = 3.464990813210079
$ws.Range("B19").Value = 0.8014218085242817
$ws.Range("C19").Value = 0.160973745007162
$ws.Range("D19").Value = 0.05125880530342641
$ws.Range("E19").Value = 0.115995163672558
$ws.Range("F19").Value = 0.9835225527429401
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("K19").Value = 0.458802513509994
$ws.Range("L19").Value = 0.2178363823606446
$ws.Range("N19").Value = 1.726401020780411
$ws.Range("O19").Value = 3.465233088753081
$ws.Range("B20").Value = 0.8197119948320903
$ws.Range("C20").Value = 0.1610594922270039
$ws.Range("D20").Value = 0.05200980928197652
$ws.Range("E20").Value = 0.1163571262242193
$ws.Range("F20").Value = 0.9849791212726018
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("K20").Value = 0.4742242261522733
$ws.Range("L20").Value = 0.2210344501058614
$ws.Range("N20").Value = 1.720008370237496
$ws.Range("O20").Value = 3.464261259627619
$ws.Range("B21").Value = 0.8814395842551903
$ws.Range("C21").Value = 0.1613592449455865
$ws.Range("D21").Value = 0.05451691605150444
$ws.Range("E21").Value = 0.1176400788153131
$ws.Range("F21").Value = 0.9904988254851617
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("K21").Value = 0.5260752291709139
$ws.Range("L21").Value = 0.2319013928573241
$ws.Range("N21").Value = 1.699201417196032
$ws.Range("O21").Value = 3.463068457935236
$ws.Range("B22").Value = 0.9219835614671865
$ws.Range("C22").Value = 0.1615637341537735
$ws.Range("D22").Value = 0.05614328684676195
$ws.Range("E22").Value = 0.1185282570389283
$ws.Range("F22").Value = 0.99457244466943
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("K22").Value = 0.5599870924887966
$ws.Range("L22").Value = 0.23909398894547
$ws.Range("N22").Value = 1.686103658704935
$ws.Range("O22").Value = 3.463833171941815
$ws.Range("B23").Value = 0.9003262772671405
$ws.Range("C23").Value = 0.1614538155716403
$ws.Range("D23").Value = 0.05527638335660612
$ws.Range("E23").Value = 0.1180496788865746
$ws.Range("F23").Value = 0.9923556399579923
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("K23").Value = 0.5418856497702791
$ws.Range("L23").Value = 0.2352469181525549
$ws.Range("N23").Value = 1.693048335354863
$ws.Range("O23").Value = 3.463283748597632
$ws.Range("B24").Value = 0.8187065412446373
$ws.Range("C24").Value = 0.1610547392956505
$ws.Range("D24").Value = 0.05196862823083848
$ws.Range("E24").Value = 0.1163369973865755
$ws.Range("F24").Value = 0.9848967784005822
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("K24").Value = 0.4733771964214952
$ws.Range("L24").Value = 0.2208583679664571
$ws.Range("N24").Value = 1.720356898755035
$ws.Range("O24").Value = 3.464306830548281
$ws.Range("B25").Value = 0.7315541259502822
$ws.Range("C25").Value = 0.1606617853958952
$ws.Range("D25").Value = 0.04834920894109729
$ws.Range("E25").Value = 0.1147035761416006
$ws.Range("F25").Value = 0.9788543130697605
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("K25").Value = 0.3996011089199101
$ws.Range("L25").Value = 0.2057288020057086
$ws.Range("N25").Value = 1.751961663483052
$ws.Range("O25").Value = 3.47204249641868
